# Safer logging to avoid leaking API keys
# - Update "Proof of Concept" success criterion (G3) to also mention Recall
# - Loosen the ROI efficiency threshold (G5) from 0.01 to 0.05
# - Move the active selection to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "ROI (=ΔPerformance​/Δcost): Efficiency Check. Does the same model perform better when focused? Threshold: ΔF1>0.05 for 5x the cost (since 5 labels are used)."
$ws.Range("G3").Value = "Proof of Concept: F1, Recall > Baseline. (Proves examples help)."

[void]$ws.Range("G3").Select()
